$wb = $excel.ActiveWorkbook
$wsStd = $wb.Worksheets.Item("Stundenerfassung")
$wsProj = $wb.Worksheets.Item("Projekte")

# Row 27 previously held "Design View Model" / "Erstellung Konzept".
# It now becomes "Erstellung ETIC2" / "Design View Model".
$wsStd.Range("B27").Value = "Erstellung ETIC2"
$wsStd.Range("C27").Value = "Design View Model"

# New row 28: 2017-06-05, "Erstellung ETIC2" / "Codierung nach MVVM", 2 hours.
$wsStd.Range("A28").Value = 42891
$wsStd.Range("A27").Copy() | Out-Null
$wsStd.Range("A28").PasteSpecial(-4122) | Out-Null
$wsStd.Range("B28").Value = "Erstellung ETIC2"
$wsStd.Range("C28").Value = "Codierung nach MVVM"
$wsStd.Range("D28").Value = 2

# Update the recorded selections on both sheets. "Stundenerfassung" must stay
# the active/tab-selected sheet, so select it last.
$wsProj.Range("B6").Select() | Out-Null
$wsStd.Range("C31").Select() | Out-Null
